$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:AF2").ClearContents()
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.59
$ws.Range("AI2").ClearContents()
$ws.Range("AJ2").Value = 23185880

# Row 3
$ws.Range("D3").Value = 1435
$ws.Range("E3").Value = 297
$ws.Range("F3").Value = 297
$ws.Range("G3").Value = 233
$ws.Range("H3").Value = 231
$ws.Range("I3").Value = 231
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 5457
$ws.Range("L3").Value = 3336
$ws.Range("M3").Value = 2121
$ws.Range("N3").Value = 2121
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 58
$ws.Range("Q3").Value = 540
$ws.Range("R3").Value = -418
$ws.Range("S3").Value = -127
$ws.Range("T3").Value = 432
$ws.Range("U3").Value = 108
$ws.Range("V3").Value = 3059
$ws.Range("W3").Value = 20.72
$ws.Range("X3").Value = 16.13
$ws.Range("Y3:Z3").ClearContents()
$ws.Range("AA3").Value = 157.24
$ws.Range("AB3").Value = 3355.22
$ws.Range("AC3").Value = 998
$ws.Range("AD3").Value = 8.619999999999999
$ws.Range("AE3").Value = 9329
$ws.Range("AF3").Value = 0.92
$ws.Range("AG3").Value = 130
$ws.Range("AH3").Value = 1.51
$ws.Range("AI3").Value = 12.53
$ws.Range("AJ3").Value = 23185880

# Row 4
$ws.Range("D4").Value = 1411
$ws.Range("E4").Value = 326
$ws.Range("F4").Value = 326
$ws.Range("G4").Value = 134
$ws.Range("H4").Value = 133
$ws.Range("I4").Value = 132
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6480
$ws.Range("L4").Value = 4120
$ws.Range("M4").Value = 2360
$ws.Range("N4").Value = 2360
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 116
$ws.Range("Q4").Value = 550
$ws.Range("R4").Value = -269
$ws.Range("S4").Value = -196
$ws.Range("T4").Value = 268
$ws.Range("U4").Value = 282
$ws.Range("V4").Value = 3870
$ws.Range("W4").Value = 23.11
$ws.Range("X4").Value = 9.4
$ws.Range("Y4").Value = 5.91
$ws.Range("Z4").Value = 2.22
$ws.Range("AA4").Value = 174.56
$ws.Range("AB4").Value = 1722.55
$ws.Range("AC4").Value = 571
$ws.Range("AD4").Value = 14.17
$ws.Range("AE4").Value = 10265
$ws.Range("AF4").Value = 0.79
$ws.Range("AG4").Value = 170
$ws.Range("AH4").Value = 2.1
$ws.Range("AI4").Value = 29.55
$ws.Range("AJ4").Value = 23185880

# Row 5
$ws.Range("D5").Value = 1776
$ws.Range("E5").Value = 451
$ws.Range("F5").Value = 451
$ws.Range("G5").Value = 257
$ws.Range("H5").Value = 257
$ws.Range("I5").Value = 257
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 9048
$ws.Range("L5").Value = 6774
$ws.Range("M5").Value = 2274
$ws.Range("N5").Value = 2274
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 116
$ws.Range("Q5").Value = 761
$ws.Range("R5").Value = -29
$ws.Range("S5").Value = -394
$ws.Range("T5").Value = 131
$ws.Range("U5").Value = 631
$ws.Range("V5").Value = 6479
$ws.Range("W5").Value = 25.39
$ws.Range("X5").Value = 14.47
$ws.Range("Y5").Value = 11.09
$ws.Range("Z5").Value = 3.31
$ws.Range("AA5").Value = 297.9
$ws.Range("AB5").Value = 1900.8
$ws.Range("AC5").Value = 1108
$ws.Range("AD5").Value = 8.19
$ws.Range("AE5").Value = 9990
$ws.Range("AF5").Value = 0.91
$ws.Range("AG5").Value = 230
$ws.Range("AH5").Value = 2.54
$ws.Range("AI5").Value = 20.38
$ws.Range("AJ5").Value = 23185880

# Row 6
$ws.Range("D6").Value = 2025
$ws.Range("E6").Value = 471
$ws.Range("F6").Value = 471
$ws.Range("G6").Value = 210
$ws.Range("H6").Value = 208
$ws.Range("I6").Value = 208
$ws.Range("K6").Value = 9026
$ws.Range("L6").Value = 6532
$ws.Range("M6").Value = 2494
$ws.Range("N6").Value = 2494
$ws.Range("P6").Value = 116
$ws.Range("Q6").Value = 825
$ws.Range("R6").Value = 12
$ws.Range("S6").Value = -1079
$ws.Range("T6").Value = 49
$ws.Range("U6").Value = 776
$ws.Range("V6").Value = 6136
$ws.Range("W6").Value = 23.26
$ws.Range("X6").Value = 10.28
$ws.Range("Y6").Value = 8.73
$ws.Range("Z6").Value = 2.3
$ws.Range("AA6").Value = 261.93
$ws.Range("AB6").Value = 2024.62
$ws.Range("AC6").Value = 898
$ws.Range("AD6").Value = 7.22
$ws.Range("AE6").Value = 11140
$ws.Range("AF6").Value = 0.58
$ws.Range("AG6:AH6").ClearContents()
$ws.Range("AI6").Value = 24.75
$ws.Range("AJ6").Value = 23185880

# Row 7
$ws.Range("D7").Value = 2273
$ws.Range("E7").Value = 492
$ws.Range("G7").Value = 171
$ws.Range("H7").Value = 169
$ws.Range("I7").Value = 169
$ws.Range("K7").Value = 8583
$ws.Range("L7").Value = 5929
$ws.Range("M7").Value = 2654
$ws.Range("N7").Value = 2654
$ws.Range("P7").Value = 116
$ws.Range("Q7").Value = 490
$ws.Range("R7").Value = -46
$ws.Range("S7").Value = -1926
$ws.Range("T7").Value = 45
$ws.Range("U7").ClearContents()
$ws.Range("W7").Value = 21.64
$ws.Range("X7").Value = 7.43
$ws.Range("Y7").Value = 6.57
$ws.Range("Z7").Value = 1.92
$ws.Range("AA7").Value = 223.4
$ws.Range("AC7").Value = 729
$ws.Range("AD7").Value = 9.58
$ws.Range("AE7").Value = 11817
$ws.Range("AF7").Value = 0.59
$ws.Range("AG7").Value = 230
$ws.Range("AH7").Value = 3.3
$ws.Range("AI7").Value = 31.55

# Row 8
$ws.Range("D8").Value = 2305
$ws.Range("E8").Value = 506
$ws.Range("G8").Value = 253
$ws.Range("H8").Value = 251
$ws.Range("I8").Value = 251
$ws.Range("K8").Value = 8278
$ws.Range("L8").Value = 5354
$ws.Range("M8").Value = 2924
$ws.Range("N8").Value = 2924
$ws.Range("P8").Value = 116
$ws.Range("Q8").Value = 580
$ws.Range("R8").Value = -68
$ws.Range("S8").Value = -1292
$ws.Range("T8").Value = 67
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 21.95
$ws.Range("X8").Value = 10.89
$ws.Range("Y8").Value = 9
$ws.Range("Z8").Value = 2.98
$ws.Range("AA8").Value = 183.11
$ws.Range("AC8").Value = 1083
$ws.Range("AD8").Value = 6.45
$ws.Range("AE8").Value = 13019
$ws.Range("AF8").Value = 0.54
$ws.Range("AG8").Value = 230
$ws.Range("AH8").Value = 3.3
$ws.Range("AI8").Value = 21.25

# Row 9
$ws.Range("D9").Value = 2852
$ws.Range("E9").Value = 666
$ws.Range("G9").Value = 411
$ws.Range("H9").Value = 407
$ws.Range("I9").Value = 407
$ws.Range("K9").Value = 9015
$ws.Range("L9").Value = 5682
$ws.Range("M9").Value = 3333
$ws.Range("N9").Value = 3333
$ws.Range("P9").Value = 116
$ws.Range("Q9").Value = 697
$ws.Range("R9").Value = -651
$ws.Range("S9").Value = -420
$ws.Range("T9").Value = 650
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 23.35
$ws.Range("X9").Value = 14.27
$ws.Range("Y9").Value = 13.01
$ws.Range("Z9").Value = 4.71
$ws.Range("AA9").Value = 170.48
$ws.Range("AC9").Value = 1755
$ws.Range("AD9").Value = 3.98
$ws.Range("AE9").Value = 14840
$ws.Range("AF9").Value = 0.47
$ws.Range("AG9").Value = 250
$ws.Range("AH9").Value = 3.58
$ws.Range("AI9").Value = 14.24
